# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-20 09:14:51
#
# For every row in the "Recorded By" column (G), when the value is a
# comma-separated list of recorders (e.g. "System, someone@example.com"),
# swap the first and last entries in that list (e.g.
# "System, someone@example.com" -> "someone@example.com, System").
# Single-value cells (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value
    $parts = $text -split ', '

    if ($parts.Count -ge 2) {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $newText = [string]::Join(', ', $parts)
        $cell.Value = $newText
    }
}
